$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D1").Value = "Memory Usage (bytes)"

$ws.Range("C2").Value = 16.50881767272949
$ws.Range("C3").Value = 15.82217216491699
$ws.Range("C4").Value = 15.93804359436035
$ws.Range("C5").Value = 16.19291305541992
$ws.Range("C6").Value = 16.74079895019531
